# Update "horarios" workbook (Línea 141) with the newest scrape (07:38:09).
#
# Sheet "LP1912": the scrap at 07:13:03 produced rows 40-42; the new scrap at
# 07:38:09 produced one more "23_HERNANDEZ" arrival (now the first of this
# batch) plus 5 brand-new arrivals. So row 40 becomes the new scrap's first
# row, the old rows 40-42 shift down to 41-43, and rows 44-48 are appended.

$wb = $excel.ActiveWorkbook

$oldStamp = "07:13:03"
$newStamp = "07:38:09"

# ---- Sheet: LP1912 --------------------------------------------------
$ws = $wb.Worksheets.Item("LP1912")

$ws.Range("A2").Value = "Última actualización: $newStamp"
$ws.Range("A3").Value = "Total filas: 43"

# Shift former rows 42, 41, 40 down to 43, 42, 41 (process bottom-up so we
# don't clobber data before it's copied).
$ws.Cells.Item(43, 1).Value = $oldStamp
$ws.Cells.Item(43, 2).Value = "08:58"
$ws.Cells.Item(43, 3).Value = "215A_EL PATO"
$ws.Cells.Item(43, 4).Value = 105
$ws.Cells.Item(43, 5).Value = "LP1912"

$ws.Cells.Item(42, 1).Value = $oldStamp
$ws.Cells.Item(42, 2).Value = "08:54"
$ws.Cells.Item(42, 3).Value = "215B_EL PATO"
$ws.Cells.Item(42, 4).Value = 101
$ws.Cells.Item(42, 5).Value = "LP1912"

$ws.Cells.Item(41, 1).Value = $oldStamp
$ws.Cells.Item(41, 2).Value = "08:52"
$ws.Cells.Item(41, 3).Value = "23_HERNANDEZ"
$ws.Cells.Item(41, 4).Value = 99
$ws.Cells.Item(41, 5).Value = "LP1912"

# Row 40 now holds the first entry of the new 07:38:09 scrap.
$ws.Cells.Item(40, 1).Value = $newStamp
$ws.Cells.Item(40, 2).Value = "08:51"
$ws.Cells.Item(40, 3).Value = "23_HERNANDEZ"
$ws.Cells.Item(40, 4).Value = 73
$ws.Cells.Item(40, 5).Value = "LP1912"

# New rows appended from the 07:38:09 scrap.
$ws.Cells.Item(44, 1).Value = $newStamp
$ws.Cells.Item(44, 2).Value = "09:06"
$ws.Cells.Item(44, 3).Value = "16_SANTA ANA"
$ws.Cells.Item(44, 4).Value = 88
$ws.Cells.Item(44, 5).Value = "LP1912"

$ws.Cells.Item(45, 1).Value = $newStamp
$ws.Cells.Item(45, 2).Value = "09:14"
$ws.Cells.Item(45, 3).Value = "27_EL RETIRO"
$ws.Cells.Item(45, 4).Value = 96
$ws.Cells.Item(45, 5).Value = "LP1912"

$ws.Cells.Item(46, 1).Value = $newStamp
$ws.Cells.Item(46, 2).Value = "09:18"
$ws.Cells.Item(46, 3).Value = "15X38_ABASTO"
$ws.Cells.Item(46, 4).Value = 100
$ws.Cells.Item(46, 5).Value = "LP1912"

$ws.Cells.Item(47, 1).Value = $newStamp
$ws.Cells.Item(47, 2).Value = "09:18"
$ws.Cells.Item(47, 3).Value = "14_ABASTO"
$ws.Cells.Item(47, 4).Value = 100
$ws.Cells.Item(47, 5).Value = "LP1912"

$ws.Cells.Item(48, 1).Value = $newStamp
$ws.Cells.Item(48, 2).Value = "09:29"
$ws.Cells.Item(48, 3).Value = "10_OLMOS"
$ws.Cells.Item(48, 4).Value = 111
$ws.Cells.Item(48, 5).Value = "LP1912"

# ---- Sheet: LP1912-215 -----------------------------------------------
$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws2.Range("A2").Value = "Última actualización: $newStamp"

# ---- Sheet: 6203-6173 -------------------------------------------------
$ws3 = $wb.Worksheets.Item("6203-6173")
$ws3.Range("A2").Value = "Última actualización: $newStamp"
